$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row content permutation (records physically re-sorted within the sheet) ---
# Row 5 <- source data from original row 6
$ws.Range("A5:AY5").ClearContents()
$ws.Range("A5").Value2 = 131256691
$ws.Range("B5").Value2 = 57884
$ws.Range("D5").Value2 = 'NT'
$ws.Range("E5").Value2 = 100109
$ws.Range("F5").Value2 = 'Tretåig hackspett'
$ws.Range("G5").Value2 = 'Picoides tridactylus'
$ws.Range("H5").Value2 = '(Linnaeus, 1758)'
$ws.Range("I5").Value2 = ''
$ws.Range("M5").Value2 = 'äldre spår'
$ws.Range("P5").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q5").Value2 = 488667
$ws.Range("R5").Value2 = 6665262
$ws.Range("S5").Value2 = 5
$ws.Range("T5").Value2 = 'Dalarna'
$ws.Range("U5").Value2 = 'Ludvika'
$ws.Range("V5").Value2 = 'Dalarna'
$ws.Range("W5").Value2 = 'Grangärde'
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value2 = '2026-02-22'
$ws.Range("Z5").Value2 = '10:55'
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value2 = '2026-02-22'
$ws.Range("AB5").Value2 = '10:55'
$ws.Range("AC5").Value2 = 'Ringhack på gran.'
$ws.Range("AD5").Value2 = $false
$ws.Range("AE5").Value2 = $false
$ws.Range("AG5").Value2 = $false
$ws.Range("AT5").Value2 = ''
$ws.Range("AW5").Value2 = 'Karl Ericson'
$ws.Range("AX5").Value2 = 'Karl Ericson'
$ws.Range("AY5").Value2 = ''

# Row 6 <- source data from original row 5
$ws.Range("A6:AY6").ClearContents()
$ws.Range("A6").Value2 = 131260583
$ws.Range("B6").Value2 = 57884
$ws.Range("D6").Value2 = 'NT'
$ws.Range("E6").Value2 = 100109
$ws.Range("F6").Value2 = 'Tretåig hackspett'
$ws.Range("G6").Value2 = 'Picoides tridactylus'
$ws.Range("H6").Value2 = '(Linnaeus, 1758)'
$ws.Range("I6").Value2 = ''
$ws.Range("M6").Value2 = 'färska spår'
$ws.Range("P6").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q6").Value2 = 488834
$ws.Range("R6").Value2 = 6665228
$ws.Range("S6").Value2 = 5
$ws.Range("T6").Value2 = 'Dalarna'
$ws.Range("U6").Value2 = 'Ludvika'
$ws.Range("V6").Value2 = 'Dalarna'
$ws.Range("W6").Value2 = 'Grangärde'
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value2 = '2026-02-22'
$ws.Range("Z6").Value2 = '15:30'
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value2 = '2026-02-22'
$ws.Range("AB6").Value2 = '15:30'
$ws.Range("AC6").Value2 = 'Ringhack på tall.'
$ws.Range("AD6").Value2 = $false
$ws.Range("AE6").Value2 = $false
$ws.Range("AG6").Value2 = $false
$ws.Range("AT6").Value2 = ''
$ws.Range("AW6").Value2 = 'Karl Ericson'
$ws.Range("AX6").Value2 = 'Karl Ericson'
$ws.Range("AY6").Value2 = ''

# Row 13 <- source data from original row 14
$ws.Range("A13:AY13").ClearContents()
$ws.Range("A13").Value2 = 131260641
$ws.Range("B13").Value2 = 57884
$ws.Range("D13").Value2 = 'NT'
$ws.Range("E13").Value2 = 100109
$ws.Range("F13").Value2 = 'Tretåig hackspett'
$ws.Range("G13").Value2 = 'Picoides tridactylus'
$ws.Range("H13").Value2 = '(Linnaeus, 1758)'
$ws.Range("I13").Value2 = ''
$ws.Range("M13").Value2 = 'äldre spår'
$ws.Range("P13").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q13").Value2 = 488859
$ws.Range("R13").Value2 = 6665292
$ws.Range("S13").Value2 = 5
$ws.Range("T13").Value2 = 'Dalarna'
$ws.Range("U13").Value2 = 'Ludvika'
$ws.Range("V13").Value2 = 'Dalarna'
$ws.Range("W13").Value2 = 'Grangärde'
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value2 = '2026-02-22'
$ws.Range("Z13").Value2 = '15:34'
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value2 = '2026-02-22'
$ws.Range("AB13").Value2 = '15:34'
$ws.Range("AC13").Value2 = 'Ringhack på gran.'
$ws.Range("AD13").Value2 = $false
$ws.Range("AE13").Value2 = $false
$ws.Range("AG13").Value2 = $false
$ws.Range("AT13").Value2 = ''
$ws.Range("AW13").Value2 = 'Karl Ericson'
$ws.Range("AX13").Value2 = 'Karl Ericson'
$ws.Range("AY13").Value2 = ''

# Row 14 <- source data from original row 15
$ws.Range("A14:AY14").ClearContents()
$ws.Range("A14").Value2 = 131257290
$ws.Range("B14").Value2 = 57884
$ws.Range("D14").Value2 = 'NT'
$ws.Range("E14").Value2 = 100109
$ws.Range("F14").Value2 = 'Tretåig hackspett'
$ws.Range("G14").Value2 = 'Picoides tridactylus'
$ws.Range("H14").Value2 = '(Linnaeus, 1758)'
$ws.Range("I14").Value2 = ''
$ws.Range("M14").Value2 = 'äldre spår'
$ws.Range("P14").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q14").Value2 = 488842
$ws.Range("R14").Value2 = 6665224
$ws.Range("S14").Value2 = 5
$ws.Range("T14").Value2 = 'Dalarna'
$ws.Range("U14").Value2 = 'Ludvika'
$ws.Range("V14").Value2 = 'Dalarna'
$ws.Range("W14").Value2 = 'Grangärde'
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value2 = '2026-02-22'
$ws.Range("Z14").Value2 = '11:26'
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value2 = '2026-02-22'
$ws.Range("AB14").Value2 = '11:26'
$ws.Range("AC14").Value2 = 'Ringhack på tall.'
$ws.Range("AD14").Value2 = $false
$ws.Range("AE14").Value2 = $false
$ws.Range("AG14").Value2 = $false
$ws.Range("AT14").Value2 = ''
$ws.Range("AW14").Value2 = 'Karl Ericson'
$ws.Range("AX14").Value2 = 'Karl Ericson'
$ws.Range("AY14").Value2 = ''

# Row 15 <- source data from original row 13
$ws.Range("A15:AY15").ClearContents()
$ws.Range("A15").Value2 = 131256673
$ws.Range("B15").Value2 = 57884
$ws.Range("D15").Value2 = 'NT'
$ws.Range("E15").Value2 = 100109
$ws.Range("F15").Value2 = 'Tretåig hackspett'
$ws.Range("G15").Value2 = 'Picoides tridactylus'
$ws.Range("H15").Value2 = '(Linnaeus, 1758)'
$ws.Range("I15").Value2 = ''
$ws.Range("M15").Value2 = 'äldre spår'
$ws.Range("P15").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q15").Value2 = 488652
$ws.Range("R15").Value2 = 6665282
$ws.Range("S15").Value2 = 5
$ws.Range("T15").Value2 = 'Dalarna'
$ws.Range("U15").Value2 = 'Ludvika'
$ws.Range("V15").Value2 = 'Dalarna'
$ws.Range("W15").Value2 = 'Grangärde'
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value2 = '2026-02-22'
$ws.Range("Z15").Value2 = '10:54'
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value2 = '2026-02-22'
$ws.Range("AB15").Value2 = '10:54'
$ws.Range("AC15").Value2 = 'Ringhack på tall.'
$ws.Range("AD15").Value2 = $false
$ws.Range("AE15").Value2 = $false
$ws.Range("AG15").Value2 = $false
$ws.Range("AT15").Value2 = ''
$ws.Range("AW15").Value2 = 'Karl Ericson'
$ws.Range("AX15").Value2 = 'Karl Ericson'
$ws.Range("AY15").Value2 = ''

# Row 25 <- source data from original row 26
$ws.Range("A25:AY25").ClearContents()
$ws.Range("A25").Value2 = 131257045
$ws.Range("B25").Value2 = 79244
$ws.Range("D25").Value2 = 'NT'
$ws.Range("E25").Value2 = 6425
$ws.Range("F25").Value2 = 'Garnlav'
$ws.Range("G25").Value2 = 'Alectoria sarmentosa'
$ws.Range("H25").Value2 = '(Ach.) Ach.'
$ws.Range("I25").Value2 = ''
$ws.Range("P25").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q25").Value2 = 488760
$ws.Range("R25").Value2 = 6665302
$ws.Range("S25").Value2 = 5
$ws.Range("T25").Value2 = 'Dalarna'
$ws.Range("U25").Value2 = 'Ludvika'
$ws.Range("V25").Value2 = 'Dalarna'
$ws.Range("W25").Value2 = 'Grangärde'
$ws.Range("Y25").NumberFormat = "@"
$ws.Range("Y25").Value2 = '2026-02-22'
$ws.Range("Z25").Value2 = '11:10'
$ws.Range("AA25").NumberFormat = "@"
$ws.Range("AA25").Value2 = '2026-02-22'
$ws.Range("AB25").Value2 = '11:10'
$ws.Range("AC25").Value2 = 'Gran.'
$ws.Range("AD25").Value2 = $false
$ws.Range("AE25").Value2 = $false
$ws.Range("AG25").Value2 = $false
$ws.Range("AT25").Value2 = ''
$ws.Range("AW25").Value2 = 'Karl Ericson'
$ws.Range("AX25").Value2 = 'Karl Ericson'
$ws.Range("AY25").Value2 = ''

# Row 26 <- source data from original row 25
$ws.Range("A26:AY26").ClearContents()
$ws.Range("A26").Value2 = 131257650
$ws.Range("B26").Value2 = 79244
$ws.Range("D26").Value2 = 'NT'
$ws.Range("E26").Value2 = 6425
$ws.Range("F26").Value2 = 'Garnlav'
$ws.Range("G26").Value2 = 'Alectoria sarmentosa'
$ws.Range("H26").Value2 = '(Ach.) Ach.'
$ws.Range("I26").Value2 = ''
$ws.Range("P26").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q26").Value2 = 488911
$ws.Range("R26").Value2 = 6665227
$ws.Range("S26").Value2 = 5
$ws.Range("T26").Value2 = 'Dalarna'
$ws.Range("U26").Value2 = 'Ludvika'
$ws.Range("V26").Value2 = 'Dalarna'
$ws.Range("W26").Value2 = 'Grangärde'
$ws.Range("Y26").NumberFormat = "@"
$ws.Range("Y26").Value2 = '2026-02-22'
$ws.Range("Z26").Value2 = '12:00'
$ws.Range("AA26").NumberFormat = "@"
$ws.Range("AA26").Value2 = '2026-02-22'
$ws.Range("AB26").Value2 = '12:00'
$ws.Range("AC26").Value2 = 'Gran'
$ws.Range("AD26").Value2 = $false
$ws.Range("AE26").Value2 = $false
$ws.Range("AG26").Value2 = $false
$ws.Range("AT26").Value2 = ''
$ws.Range("AW26").Value2 = 'Karl Ericson'
$ws.Range("AX26").Value2 = 'Karl Ericson'
$ws.Range("AY26").Value2 = ''

# Row 28 <- source data from original row 29
$ws.Range("A28:AY28").ClearContents()
$ws.Range("A28").Value2 = 131256649
$ws.Range("B28").Value2 = 79244
$ws.Range("D28").Value2 = 'NT'
$ws.Range("E28").Value2 = 6425
$ws.Range("F28").Value2 = 'Garnlav'
$ws.Range("G28").Value2 = 'Alectoria sarmentosa'
$ws.Range("H28").Value2 = '(Ach.) Ach.'
$ws.Range("I28").Value2 = ''
$ws.Range("P28").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q28").Value2 = 488685
$ws.Range("R28").Value2 = 6665288
$ws.Range("S28").Value2 = 5
$ws.Range("T28").Value2 = 'Dalarna'
$ws.Range("U28").Value2 = 'Ludvika'
$ws.Range("V28").Value2 = 'Dalarna'
$ws.Range("W28").Value2 = 'Grangärde'
$ws.Range("Y28").NumberFormat = "@"
$ws.Range("Y28").Value2 = '2026-02-22'
$ws.Range("Z28").Value2 = '10:52'
$ws.Range("AA28").NumberFormat = "@"
$ws.Range("AA28").Value2 = '2026-02-22'
$ws.Range("AB28").Value2 = '10:52'
$ws.Range("AC28").Value2 = 'Gran.'
$ws.Range("AD28").Value2 = $false
$ws.Range("AE28").Value2 = $false
$ws.Range("AG28").Value2 = $false
$ws.Range("AT28").Value2 = ''
$ws.Range("AW28").Value2 = 'Karl Ericson'
$ws.Range("AX28").Value2 = 'Karl Ericson'
$ws.Range("AY28").Value2 = ''

# Row 29 <- source data from original row 30
$ws.Range("A29:AY29").ClearContents()
$ws.Range("A29").Value2 = 131255910
$ws.Range("B29").Value2 = 79244
$ws.Range("D29").Value2 = 'NT'
$ws.Range("E29").Value2 = 6425
$ws.Range("F29").Value2 = 'Garnlav'
$ws.Range("G29").Value2 = 'Alectoria sarmentosa'
$ws.Range("H29").Value2 = '(Ach.) Ach.'
$ws.Range("I29").Value2 = ''
$ws.Range("P29").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q29").Value2 = 488763
$ws.Range("R29").Value2 = 6665157
$ws.Range("S29").Value2 = 5
$ws.Range("T29").Value2 = 'Dalarna'
$ws.Range("U29").Value2 = 'Ludvika'
$ws.Range("V29").Value2 = 'Dalarna'
$ws.Range("W29").Value2 = 'Grangärde'
$ws.Range("Y29").NumberFormat = "@"
$ws.Range("Y29").Value2 = '2026-02-22'
$ws.Range("Z29").Value2 = '10:03'
$ws.Range("AA29").NumberFormat = "@"
$ws.Range("AA29").Value2 = '2026-02-22'
$ws.Range("AB29").Value2 = '10:03'
$ws.Range("AC29").Value2 = 'Tall.'
$ws.Range("AD29").Value2 = $false
$ws.Range("AE29").Value2 = $false
$ws.Range("AG29").Value2 = $false
$ws.Range("AT29").Value2 = ''
$ws.Range("AW29").Value2 = 'Karl Ericson'
$ws.Range("AX29").Value2 = 'Karl Ericson'
$ws.Range("AY29").Value2 = ''

# Row 30 <- source data from original row 31
$ws.Range("A30:AY30").ClearContents()
$ws.Range("A30").Value2 = 131258531
$ws.Range("B30").Value2 = 79244
$ws.Range("D30").Value2 = 'NT'
$ws.Range("E30").Value2 = 6425
$ws.Range("F30").Value2 = 'Garnlav'
$ws.Range("G30").Value2 = 'Alectoria sarmentosa'
$ws.Range("H30").Value2 = '(Ach.) Ach.'
$ws.Range("I30").Value2 = ''
$ws.Range("P30").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q30").Value2 = 488725
$ws.Range("R30").Value2 = 6665212
$ws.Range("S30").Value2 = 5
$ws.Range("T30").Value2 = 'Dalarna'
$ws.Range("U30").Value2 = 'Ludvika'
$ws.Range("V30").Value2 = 'Dalarna'
$ws.Range("W30").Value2 = 'Grangärde'
$ws.Range("Y30").NumberFormat = "@"
$ws.Range("Y30").Value2 = '2026-02-22'
$ws.Range("Z30").Value2 = '13:02'
$ws.Range("AA30").NumberFormat = "@"
$ws.Range("AA30").Value2 = '2026-02-22'
$ws.Range("AB30").Value2 = '13:02'
$ws.Range("AC30").Value2 = 'Gran'
$ws.Range("AD30").Value2 = $false
$ws.Range("AE30").Value2 = $false
$ws.Range("AG30").Value2 = $false
$ws.Range("AT30").Value2 = ''
$ws.Range("AW30").Value2 = 'Karl Ericson'
$ws.Range("AX30").Value2 = 'Karl Ericson'
$ws.Range("AY30").Value2 = ''

# Row 31 <- source data from original row 28
$ws.Range("A31:AY31").ClearContents()
$ws.Range("A31").Value2 = 131257239
$ws.Range("B31").Value2 = 57884
$ws.Range("D31").Value2 = 'NT'
$ws.Range("E31").Value2 = 100109
$ws.Range("F31").Value2 = 'Tretåig hackspett'
$ws.Range("G31").Value2 = 'Picoides tridactylus'
$ws.Range("H31").Value2 = '(Linnaeus, 1758)'
$ws.Range("I31").Value2 = ''
$ws.Range("M31").Value2 = 'färska spår'
$ws.Range("P31").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q31").Value2 = 488852
$ws.Range("R31").Value2 = 6665286
$ws.Range("S31").Value2 = 5
$ws.Range("T31").Value2 = 'Dalarna'
$ws.Range("U31").Value2 = 'Ludvika'
$ws.Range("V31").Value2 = 'Dalarna'
$ws.Range("W31").Value2 = 'Grangärde'
$ws.Range("Y31").NumberFormat = "@"
$ws.Range("Y31").Value2 = '2026-02-22'
$ws.Range("Z31").Value2 = '11:23'
$ws.Range("AA31").NumberFormat = "@"
$ws.Range("AA31").Value2 = '2026-02-22'
$ws.Range("AB31").Value2 = '11:23'
$ws.Range("AC31").Value2 = 'Barkfläk, hagelsalva.'
$ws.Range("AD31").Value2 = $false
$ws.Range("AE31").Value2 = $false
$ws.Range("AG31").Value2 = $false
$ws.Range("AT31").Value2 = ''
$ws.Range("AW31").Value2 = 'Karl Ericson'
$ws.Range("AX31").Value2 = 'Karl Ericson'
$ws.Range("AY31").Value2 = ''

# Row 32 <- source data from original row 34
$ws.Range("A32:AY32").ClearContents()
$ws.Range("A32").Value2 = 131258537
$ws.Range("B32").Value2 = 79244
$ws.Range("D32").Value2 = 'NT'
$ws.Range("E32").Value2 = 6425
$ws.Range("F32").Value2 = 'Garnlav'
$ws.Range("G32").Value2 = 'Alectoria sarmentosa'
$ws.Range("H32").Value2 = '(Ach.) Ach.'
$ws.Range("I32").Value2 = ''
$ws.Range("P32").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q32").Value2 = 488726
$ws.Range("R32").Value2 = 6665209
$ws.Range("S32").Value2 = 5
$ws.Range("T32").Value2 = 'Dalarna'
$ws.Range("U32").Value2 = 'Ludvika'
$ws.Range("V32").Value2 = 'Dalarna'
$ws.Range("W32").Value2 = 'Grangärde'
$ws.Range("Y32").NumberFormat = "@"
$ws.Range("Y32").Value2 = '2026-02-22'
$ws.Range("Z32").Value2 = '13:02'
$ws.Range("AA32").NumberFormat = "@"
$ws.Range("AA32").Value2 = '2026-02-22'
$ws.Range("AB32").Value2 = '13:02'
$ws.Range("AC32").Value2 = 'Gran'
$ws.Range("AD32").Value2 = $false
$ws.Range("AE32").Value2 = $false
$ws.Range("AG32").Value2 = $false
$ws.Range("AT32").Value2 = ''
$ws.Range("AW32").Value2 = 'Karl Ericson'
$ws.Range("AX32").Value2 = 'Karl Ericson'
$ws.Range("AY32").Value2 = ''

# Row 34 <- source data from original row 32
$ws.Range("A34:AY34").ClearContents()
$ws.Range("A34").Value2 = 131257659
$ws.Range("B34").Value2 = 57884
$ws.Range("D34").Value2 = 'NT'
$ws.Range("E34").Value2 = 100109
$ws.Range("F34").Value2 = 'Tretåig hackspett'
$ws.Range("G34").Value2 = 'Picoides tridactylus'
$ws.Range("H34").Value2 = '(Linnaeus, 1758)'
$ws.Range("I34").Value2 = ''
$ws.Range("M34").Value2 = 'äldre spår'
$ws.Range("P34").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q34").Value2 = 488901
$ws.Range("R34").Value2 = 6665222
$ws.Range("S34").Value2 = 5
$ws.Range("T34").Value2 = 'Dalarna'
$ws.Range("U34").Value2 = 'Ludvika'
$ws.Range("V34").Value2 = 'Dalarna'
$ws.Range("W34").Value2 = 'Grangärde'
$ws.Range("Y34").NumberFormat = "@"
$ws.Range("Y34").Value2 = '2026-02-22'
$ws.Range("Z34").Value2 = '12:02'
$ws.Range("AA34").NumberFormat = "@"
$ws.Range("AA34").Value2 = '2026-02-22'
$ws.Range("AB34").Value2 = '12:02'
$ws.Range("AC34").Value2 = 'Ringhack på tall.'
$ws.Range("AD34").Value2 = $false
$ws.Range("AE34").Value2 = $false
$ws.Range("AG34").Value2 = $false
$ws.Range("AT34").Value2 = ''
$ws.Range("AW34").Value2 = 'Karl Ericson'
$ws.Range("AX34").Value2 = 'Karl Ericson'
$ws.Range("AY34").Value2 = ''

# Row 36 <- source data from original row 37
$ws.Range("A36:AY36").ClearContents()
$ws.Range("A36").Value2 = 131257385
$ws.Range("B36").Value2 = 91829
$ws.Range("D36").Value2 = 'NT'
$ws.Range("E36").Value2 = 5432
$ws.Range("F36").Value2 = 'Granticka'
$ws.Range("G36").Value2 = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H36").Value2 = ''
$ws.Range("I36").Value2 = ''
$ws.Range("P36").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q36").Value2 = 488876
$ws.Range("R36").Value2 = 6665194
$ws.Range("S36").Value2 = 5
$ws.Range("T36").Value2 = 'Dalarna'
$ws.Range("U36").Value2 = 'Ludvika'
$ws.Range("V36").Value2 = 'Dalarna'
$ws.Range("W36").Value2 = 'Grangärde'
$ws.Range("Y36").NumberFormat = "@"
$ws.Range("Y36").Value2 = '2026-02-22'
$ws.Range("Z36").Value2 = '11:31'
$ws.Range("AA36").NumberFormat = "@"
$ws.Range("AA36").Value2 = '2026-02-22'
$ws.Range("AB36").Value2 = '11:31'
$ws.Range("AC36").Value2 = 'Lågstubbe.'
$ws.Range("AD36").Value2 = $false
$ws.Range("AE36").Value2 = $false
$ws.Range("AG36").Value2 = $false
$ws.Range("AT36").Value2 = ''
$ws.Range("AW36").Value2 = 'Karl Ericson'
$ws.Range("AX36").Value2 = 'Karl Ericson'
$ws.Range("AY36").Value2 = ''

# Row 37 <- source data from original row 38
$ws.Range("A37:AY37").ClearContents()
$ws.Range("A37").Value2 = 131260531
$ws.Range("B37").Value2 = 79244
$ws.Range("D37").Value2 = 'NT'
$ws.Range("E37").Value2 = 6425
$ws.Range("F37").Value2 = 'Garnlav'
$ws.Range("G37").Value2 = 'Alectoria sarmentosa'
$ws.Range("H37").Value2 = '(Ach.) Ach.'
$ws.Range("I37").Value2 = ''
$ws.Range("P37").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q37").Value2 = 488786
$ws.Range("R37").Value2 = 6665188
$ws.Range("S37").Value2 = 5
$ws.Range("T37").Value2 = 'Dalarna'
$ws.Range("U37").Value2 = 'Ludvika'
$ws.Range("V37").Value2 = 'Dalarna'
$ws.Range("W37").Value2 = 'Grangärde'
$ws.Range("Y37").NumberFormat = "@"
$ws.Range("Y37").Value2 = '2026-02-22'
$ws.Range("Z37").Value2 = '15:25'
$ws.Range("AA37").NumberFormat = "@"
$ws.Range("AA37").Value2 = '2026-02-22'
$ws.Range("AB37").Value2 = '15:25'
$ws.Range("AC37").Value2 = 'Gran'
$ws.Range("AD37").Value2 = $false
$ws.Range("AE37").Value2 = $false
$ws.Range("AG37").Value2 = $false
$ws.Range("AT37").Value2 = ''
$ws.Range("AW37").Value2 = 'Karl Ericson'
$ws.Range("AX37").Value2 = 'Karl Ericson'
$ws.Range("AY37").Value2 = ''

# Row 38 <- source data from original row 36
$ws.Range("A38:AY38").ClearContents()
$ws.Range("A38").Value2 = 131256459
$ws.Range("B38").Value2 = 57881
$ws.Range("D38").Value2 = 'NT'
$ws.Range("E38").Value2 = 100049
$ws.Range("F38").Value2 = 'Spillkråka'
$ws.Range("G38").Value2 = 'Dryocopus martius'
$ws.Range("H38").Value2 = '(Linnaeus, 1758)'
$ws.Range("I38").Value2 = ''
$ws.Range("M38").Value2 = 'färska spår'
$ws.Range("P38").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q38").Value2 = 488669
$ws.Range("R38").Value2 = 6665268
$ws.Range("S38").Value2 = 5
$ws.Range("T38").Value2 = 'Dalarna'
$ws.Range("U38").Value2 = 'Ludvika'
$ws.Range("V38").Value2 = 'Dalarna'
$ws.Range("W38").Value2 = 'Grangärde'
$ws.Range("Y38").NumberFormat = "@"
$ws.Range("Y38").Value2 = '2026-02-22'
$ws.Range("Z38").Value2 = '10:42'
$ws.Range("AA38").NumberFormat = "@"
$ws.Range("AA38").Value2 = '2026-02-22'
$ws.Range("AB38").Value2 = '10:42'
$ws.Range("AC38").Value2 = 'Färska och äldre hack.'
$ws.Range("AD38").Value2 = $false
$ws.Range("AE38").Value2 = $false
$ws.Range("AG38").Value2 = $false
$ws.Range("AT38").Value2 = ''
$ws.Range("AW38").Value2 = 'Karl Ericson'
$ws.Range("AX38").Value2 = 'Karl Ericson'
$ws.Range("AY38").Value2 = ''

# Row 41 <- source data from original row 42
$ws.Range("A41:AY41").ClearContents()
$ws.Range("A41").Value2 = 131257035
$ws.Range("B41").Value2 = 79244
$ws.Range("D41").Value2 = 'NT'
$ws.Range("E41").Value2 = 6425
$ws.Range("F41").Value2 = 'Garnlav'
$ws.Range("G41").Value2 = 'Alectoria sarmentosa'
$ws.Range("H41").Value2 = '(Ach.) Ach.'
$ws.Range("I41").Value2 = ''
$ws.Range("P41").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q41").Value2 = 488760
$ws.Range("R41").Value2 = 6665301
$ws.Range("S41").Value2 = 5
$ws.Range("T41").Value2 = 'Dalarna'
$ws.Range("U41").Value2 = 'Ludvika'
$ws.Range("V41").Value2 = 'Dalarna'
$ws.Range("W41").Value2 = 'Grangärde'
$ws.Range("Y41").NumberFormat = "@"
$ws.Range("Y41").Value2 = '2026-02-22'
$ws.Range("Z41").Value2 = '11:09'
$ws.Range("AA41").NumberFormat = "@"
$ws.Range("AA41").Value2 = '2026-02-22'
$ws.Range("AB41").Value2 = '11:09'
$ws.Range("AC41").Value2 = 'Gran'
$ws.Range("AD41").Value2 = $false
$ws.Range("AE41").Value2 = $false
$ws.Range("AG41").Value2 = $false
$ws.Range("AT41").Value2 = ''
$ws.Range("AW41").Value2 = 'Karl Ericson'
$ws.Range("AX41").Value2 = 'Karl Ericson'
$ws.Range("AY41").Value2 = ''

# Row 42 <- source data from original row 41
$ws.Range("A42:AY42").ClearContents()
$ws.Range("A42").Value2 = 131257343
$ws.Range("B42").Value2 = 57884
$ws.Range("D42").Value2 = 'NT'
$ws.Range("E42").Value2 = 100109
$ws.Range("F42").Value2 = 'Tretåig hackspett'
$ws.Range("G42").Value2 = 'Picoides tridactylus'
$ws.Range("H42").Value2 = '(Linnaeus, 1758)'
$ws.Range("I42").Value2 = ''
$ws.Range("M42").Value2 = 'färska spår'
$ws.Range("P42").Value2 = 'Hyttfallet, Hyttfallet, Dlr'
$ws.Range("Q42").Value2 = 488872
$ws.Range("R42").Value2 = 6665191
$ws.Range("S42").Value2 = 5
$ws.Range("T42").Value2 = 'Dalarna'
$ws.Range("U42").Value2 = 'Ludvika'
$ws.Range("V42").Value2 = 'Dalarna'
$ws.Range("W42").Value2 = 'Grangärde'
$ws.Range("Y42").NumberFormat = "@"
$ws.Range("Y42").Value2 = '2026-02-22'
$ws.Range("Z42").Value2 = '11:29'
$ws.Range("AA42").NumberFormat = "@"
$ws.Range("AA42").Value2 = '2026-02-22'
$ws.Range("AB42").Value2 = '11:29'
$ws.Range("AC42").Value2 = 'Mycket barkfläk, hagelsalvor på många träd, skalade klena senvuxna granar.'
$ws.Range("AD42").Value2 = $false
$ws.Range("AE42").Value2 = $false
$ws.Range("AG42").Value2 = $false
$ws.Range("AT42").Value2 = ''
$ws.Range("AW42").Value2 = 'Karl Ericson'
$ws.Range("AX42").Value2 = 'Karl Ericson'
$ws.Range("AY42").Value2 = ''

# Row 43 <- source data from original row 44
$ws.Range("A43:AY43").ClearContents()
$ws.Range("A43").Value2 = 131273991
$ws.Range("B43").Value2 = 79244
$ws.Range("D43").Value2 = 'NT'
$ws.Range("E43").Value2 = 6425
$ws.Range("F43").Value2 = 'Garnlav'
$ws.Range("G43").Value2 = 'Alectoria sarmentosa'
$ws.Range("H43").Value2 = '(Ach.) Ach.'
$ws.Range("I43").Value2 = ''
$ws.Range("J43").Value2 = ''
$ws.Range("K43").Value2 = ''
$ws.Range("N43").Value2 = ''
$ws.Range("P43").Value2 = 'Hyttfallet, Dlr'
$ws.Range("Q43").Value2 = 488928
$ws.Range("R43").Value2 = 6665146
$ws.Range("S43").Value2 = 50
$ws.Range("T43").Value2 = 'Dalarna'
$ws.Range("U43").Value2 = 'Ludvika'
$ws.Range("V43").Value2 = 'Dalarna'
$ws.Range("W43").Value2 = 'Grangärde'
$ws.Range("Y43").NumberFormat = "@"
$ws.Range("Y43").Value2 = '2026-02-22'
$ws.Range("AA43").NumberFormat = "@"
$ws.Range("AA43").Value2 = '2026-02-22'
$ws.Range("AD43").Value2 = $false
$ws.Range("AE43").Value2 = $false
$ws.Range("AF43").Value2 = ''
$ws.Range("AG43").Value2 = $false
$ws.Range("AT43").Value2 = ''
$ws.Range("AW43").Value2 = 'Anna-Lena Thommson'
$ws.Range("AX43").Value2 = 'Anna-Lena Thommson'
$ws.Range("AY43").Value2 = ''

# Row 44 <- source data from original row 43
$ws.Range("A44:AY44").ClearContents()
$ws.Range("A44").Value2 = 131273946
$ws.Range("B44").Value2 = 79244
$ws.Range("D44").Value2 = 'NT'
$ws.Range("E44").Value2 = 6425
$ws.Range("F44").Value2 = 'Garnlav'
$ws.Range("G44").Value2 = 'Alectoria sarmentosa'
$ws.Range("H44").Value2 = '(Ach.) Ach.'
$ws.Range("I44").Value2 = ''
$ws.Range("J44").Value2 = ''
$ws.Range("K44").Value2 = ''
$ws.Range("N44").Value2 = ''
$ws.Range("P44").Value2 = 'Hyttfallet, Dlr'
$ws.Range("Q44").Value2 = 488774
$ws.Range("R44").Value2 = 6665353
$ws.Range("S44").Value2 = 50
$ws.Range("T44").Value2 = 'Dalarna'
$ws.Range("U44").Value2 = 'Ludvika'
$ws.Range("V44").Value2 = 'Dalarna'
$ws.Range("W44").Value2 = 'Grangärde'
$ws.Range("Y44").NumberFormat = "@"
$ws.Range("Y44").Value2 = '2026-02-22'
$ws.Range("AA44").NumberFormat = "@"
$ws.Range("AA44").Value2 = '2026-02-22'
$ws.Range("AD44").Value2 = $false
$ws.Range("AE44").Value2 = $false
$ws.Range("AF44").Value2 = ''
$ws.Range("AG44").Value2 = $false
$ws.Range("AT44").Value2 = ''
$ws.Range("AW44").Value2 = 'Anna-Lena Thommson'
$ws.Range("AX44").Value2 = 'Anna-Lena Thommson'
$ws.Range("AY44").Value2 = ''

# --- Taxonsorteringsordning (column B) renumbering for unaffected rows ---
$ws.Range("B4").Value2 = 79245
$ws.Range("B7").Value2 = 79245
$ws.Range("B8").Value2 = 91830
$ws.Range("B9").Value2 = 91830
$ws.Range("B11").Value2 = 79245
$ws.Range("B12").Value2 = 79245
$ws.Range("B16").Value2 = 79245
$ws.Range("B17").Value2 = 91830
$ws.Range("B18").Value2 = 79245
$ws.Range("B19").Value2 = 79245
$ws.Range("B20").Value2 = 79245
$ws.Range("B21").Value2 = 79245
$ws.Range("B22").Value2 = 81230
$ws.Range("B23").Value2 = 79245
$ws.Range("B24").Value2 = 79245
$ws.Range("B27").Value2 = 79245
$ws.Range("B33").Value2 = 79245
$ws.Range("B35").Value2 = 79245
$ws.Range("B39").Value2 = 79245
$ws.Range("B40").Value2 = 79245
$ws.Range("B45").Value2 = 79245
